$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update result values: Exam score doubled (20 -> 40) and Total recalculated (40 -> 60)
$ws.Range("E2").Value = 40
$ws.Range("F2").Value = 60
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 60

# Update the active selection to match the saved workbook state
$ws.Range("F3").Select()
